$d = $word.ActiveDocument

# Replace all "June 12, 2022" occurrences with "June 13, 2022"
$d.Content.Find.Execute("June 12, 2022", $false, $false, $false, $false, $false,
                         $true, 1, $false, "June 13, 2022", 2)

# Replace "August 11, 2022" with "August 12, 2022"
$d.Content.Find.Execute("August 11, 2022", $false, $false, $false, $false, $false,
                         $true, 1, $false, "August 12, 2022", 2)
